$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: set all text/string values in the exact order required
#            so shared-string indices come out 18..44 matching the target workbook ----
$ws.Range("A6").Value = "SCRIPT/T01P01A/um0602.ssb"
$ws.Range("C6").Value = " Stealing Time Gears is\nunthinkable!"
$ws.Range("C7").Value = " Who would go around stealing\nthem?"
$ws.Range("D6").Value = " Это немыслимо!"
$ws.Range("D7").Value = " Кто вообще осмелился красть\nШестерни Времени?"
$ws.Range("E6").Value = " Üóï îåíúòìéíï!"
$ws.Range("E7").Value = " Ëóï âïïáþå ïòíåìéìòÿ ëñàòóû\nŠåòóåñîé Âñåíåîé?"
$ws.Range("C8").Value = " I don\'t think that stench is just\na figment of my imagination…"
$ws.Range("A8").Value = "SCRIPT/T01P01A/um0606.ssb"
$ws.Range("C9").Value = " Want proof? Fewer exploration\nteams are visiting Treasure Town lately!"
$ws.Range("C10").Value = " I bet the stench is keeping them\nall away."
$ws.Range("D8").Value = " Я не думаю, что эта вонь - \nвсего лишь плод моего воображения..."
$ws.Range("D9").Value = " Хочешь доказательство? В Город\nСокровищ стало приходить меньше команд!"
$ws.Range("D10").Value = " Полагаю, что их отпугивает\nзапах."
$ws.Range("E8").Value = " Ÿ îå äôíàý, œóï üóà âïîû - \nâòåãï ìéšû ðìïä íïåãï âïïáñàçåîéÿ..."
$ws.Range("E9").Value = " Öïœåšû äïëàèàóåìûòóâï? Â Ãïñïä\nÒïëñïâéþ òóàìï ðñéöïäéóû íåîûšå ëïíàîä!"
$ws.Range("E10").Value = " Ðïìàãàý, œóï éö ïóðôãéâàåó\nèàðàö."
$ws.Range("C11").Value = " I hear the whole guild\'s going\noff on the expedition now."
$ws.Range("C12").Value = " After the last expedition, they\ncame back with loads of treasure."
$ws.Range("C13").Value = " We expect big things from your\nupcoming expedition!"
$ws.Range("A11").Value = "SCRIPT/T01P01A/um0802.ssb"
$ws.Range("D11").Value = " Я слышал, что вся гильдия\nуходит в экспедицию."
$ws.Range("D12").Value = " В прошлый раз они принесли\nмного сокровищ."
$ws.Range("D13").Value = " Надеюсь что и эта экспедиция\nбудет такой же успешной!"
$ws.Range("E11").Value = " Ÿ òìúšàì, œóï âòÿ ãéìûäéÿ\nôöïäéó â üëòðåäéøéý."
$ws.Range("E12").Value = " Â ðñïšìúê ñàè ïîé ðñéîåòìé\níîïãï òïëñïâéþ."
$ws.Range("E13").Value = " Îàäåýòû œóï é üóà üëòðåäéøéÿ\náôäåó óàëïê çå ôòðåšîïê!"

# ---- Step 2: set the numeric "line number" values in column B ----
$ws.Range("B5").Value = 551
$ws.Range("B6").Value = 520
$ws.Range("B7").Value = 523
$ws.Range("B8").Value = 495
$ws.Range("B9").Value = 498
$ws.Range("B10").Value = 501
$ws.Range("B11").Value = 470
$ws.Range("B12").Value = 473
$ws.Range("B13").Value = 476

# ---- Step 3: row heights ----
$ws.Rows.Item(5).RowHeight = 31.8
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 21.6
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 31.2
$ws.Rows.Item(10).RowHeight = 21.6
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 21.6
$ws.Rows.Item(13).RowHeight = 31.8

# ---- Step 4: apply the thin bottom-border "group separator" style (style 6/7) ----
# to rows 5, 7, 10 across columns A:E (creates new border + cellXfs entries to match)
$rowRange = $ws.Range("A5:E5")
$rowRange.WrapText = $true
$rowRange.Borders.Item(9).LineStyle = 1
$rowRange.Borders.Item(9).Weight = 2
$ws.Range("C5:E5").Font.Size = 8
$rowRange = $ws.Range("A7:E7")
$rowRange.WrapText = $true
$rowRange.Borders.Item(9).LineStyle = 1
$rowRange.Borders.Item(9).Weight = 2
$ws.Range("C7:E7").Font.Size = 8
$rowRange = $ws.Range("A10:E10")
$rowRange.WrapText = $true
$rowRange.Borders.Item(9).LineStyle = 1
$rowRange.Borders.Item(9).Weight = 2
$ws.Range("C10:E10").Font.Size = 8

# ---- Step 5: dimension / view bookkeeping ----
$ws.Range("C13").Select()
